# "volver a generar problemas cuadraticos y lineales"
# Regenerate the random linear-problem data: update the follower
# restriction expressions/coefficients, the modified point, and the
# bf/BF/alpha vectors with newly generated values.

$wb = $excel.ActiveWorkbook

# --- Sheet: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# Row 2 (J_0_L0_v)
$ws.Range("A2").Value = "-2.0847668467443836 - 0.09833191166714061y_1 + 1.105819353615771y_2"
$ws.Range("B2").Value = "2.0847668467443836"
$ws.Range("D2").Value = "0.75"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "2.0"

# Row 3 (J_0_L0_v)
$ws.Range("A3").Value = "5.3586508387830545 - 0.17794521846270506y_1 - 1.444649796227846y_2"
$ws.Range("B3").Value = "-9.358650838783054"
$ws.Range("D3").Value = "0.19"
$ws.Range("E3").Value = "-3.5"
$ws.Range("F3").Value = "-3.1"

# Row 4 (J_0_LP_v)
$ws.Range("A4").Value = "65.31672418109892 - 2x - 4.93390077477628y_1 - 3.277425478499212y_2"
$ws.Range("B4").Value = "-81.31672418109892"
$ws.Range("D4").Value = "0.24"
$ws.Range("E4").Value = "9.8"
$ws.Range("F4").Value = "7.6"

# Row 5 (J_Ne_L0_v)
$ws.Range("A5").Value = "-61.25774950241683 + 8x + 0.6006539664486779y_1 - 0.4897640034119989y_2"
$ws.Range("B5").Value = "12.66774950241683"
$ws.Range("D5").Value = "0.92"
$ws.Range("E5").Value = "9.6"
$ws.Range("F5").Value = "0"

# Row 6 (J_Ne_L0_v)
$ws.Range("A6").Value = "-27.698080750639754 - 2x - 3.7129182068050426y_1 - 2.1007487441948633y_2"
$ws.Range("B6").Value = "-39.698080750639754"
$ws.Range("D6").Value = "0.82"
$ws.Range("E6").Value = "-3.5999999999999996"
$ws.Range("F6").Value = "-10.0"

# --- Sheet: Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = "7.1"
$ws.Range("B2").Value = "8.600000000000001"
$ws.Range("C2").Value = "2.65"

# --- Sheet: Vector_bf ---
$ws = $wb.Worksheets.Item("Vector_bf")
$ws.Range("A2").Value = "4.783685991651928"
$ws.Range("A3").Value = "2.4048979142900997"

# --- Sheet: Vector_BF ---
$ws = $wb.Worksheets.Item("Vector_BF")
$ws.Range("A2").Value = "-63.39999999999999"
$ws.Range("A3").Value = "31.596635705782617"
$ws.Range("A4").Value = "22.2015343561485"

# --- Sheet: Vector_Alpha ---
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 1.59
$ws.Range("A3").Value = 1.9500000000000002
